$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("策略更新")

# Update E5: status changed from "pending" to "passed"
$ws.Range("E5").Value = "passed"

# Update F5: note text, each strategy line now prefixed with "(ok) " (one with "(oK) ")
$newNote = "暂定手数`n(ok) dlm fl34 1`n(ok) dla fl34 2`n(ok) dljd fl34 3`n(ok) dli fl34 40`n(ok) dlv fl34 6`n(ok) dlv fl36 3`n(oK) dlm fl36 3`n(ok) dla fl36 2`n(ok) dljd fl36 1`n(ok) dli fl36 20`n(ok) dlv fw10 3`n(ok) dla fw10 2`n(ok) dla fd10 1`n(ok) dll fd10 1`n(ok) dljd fd10 1`n2。先把铁矿手数减半，我怕钱不够用`n"
$ws.Range("F5").Value = $newNote

# Update the view: last selection moved to F5, scrolled so B5 is the top-left visible cell
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("F5").Select()
